# Implement api endpoint validation
# Append one new captured-packet row to each of the four sheets.

$wb = $excel.ActiveWorkbook

# ROW50-FE-LIFTER (sheet1) -> new row 69
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$ws1.Range("A69").Value = 45758.7455387963
$ws1.Range("A69").NumberFormat = $ws1.Range("A68").NumberFormat
$ws1.Range("B69").Value = "0x01,0x90"
$ws1.Range("C69").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Range("D69").Value = "0x01,0x52"
$ws1.Range("E69").Value = "0xe"
$ws1.Range("F69").Value = 400
$ws1.Range("G69").Value = 568631262647114 * 1000000000
$ws1.Range("H69").Value = 338
$ws1.Range("I69").Value = 14

# ROW50-MID-LIFTER (sheet2) -> new row 71
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$ws2.Range("A71").Value = 45758.70621527778
$ws2.Range("A71").NumberFormat = $ws2.Range("A70").NumberFormat
$ws2.Range("B71").Value = "0x01,0x90 "
$ws2.Range("C71").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Range("D71").Value = "0x01,0x5a"
$ws2.Range("E71").Value = "0x19"
$ws2.Range("F71").Value = 400
$ws2.Range("G71").NumberFormat = "@"
$ws2.Range("G71").Value = "568631262647113771663628"
$ws2.Range("G71").ClearFormats()
$ws2.Range("H71").Value = 346
$ws2.Range("I71").Value = 25

# ROW11-FE-LIFTER (sheet3) -> new row 69
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$ws3.Range("A69").Value = 45758.77741506945
$ws3.Range("A69").NumberFormat = $ws3.Range("A68").NumberFormat
$ws3.Range("B69").Value = "0x01,0x90"
$ws3.Range("C69").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Range("D69").Value = "0x01,0x52"
$ws3.Range("E69").Value = "0x14"
$ws3.Range("F69").Value = 400
$ws3.Range("G69").Value = 568631262647114 * 1000000000
$ws3.Range("H69").Value = 338
$ws3.Range("I69").Value = 20

# ROW11-MID-LIFTER (sheet4) -> new row 69
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$ws4.Range("A69").Value = 45758.90231188657
$ws4.Range("A69").NumberFormat = $ws4.Range("A68").NumberFormat
$ws4.Range("B69").Value = "0x01,0x90"
$ws4.Range("C69").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Range("D69").Value = "0x01,0x5a"
$ws4.Range("E69").Value = "0x19"
$ws4.Range("F69").Value = 400
$ws4.Range("G69").Value = 568631262647114 * 1000000000
$ws4.Range("H69").Value = 346
$ws4.Range("I69").Value = 25

Write-Host "Rows appended."
